$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.883.50'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.084.55'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '233.23'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '0.626'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').Value = '59.35'
$ws.Range('E7').Value = '  +3.42%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.106'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').Value = '0.775'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('E15').Value = '  +2.25%  '
$ws.Range('D16').Value = '2.075.28'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('D17').Value = '37.803.90'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = '71.67'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '6.12'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '0.0₃0848'
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('D21').Value = '228.13'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').Value = '171.98'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').Value = '9.19'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('D32').Value = '4.72'
$ws.Range('E32').Value = '  +2.74%  '
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '3.42'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.82'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').Value = '5.42'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '0.0982'
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('D40').Value = '99.05'
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('D43').Value = '16.89'
$ws.Range('E43').Value = '  +8.16%  '
$ws.Range('D44').Value = '1.445.64'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '4.17'
$ws.Range('E46').Value = '  +3.01%  '
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').Value = '7.38'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '3.00'
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').Value = '2.277.08'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '46.80'
$ws.Range('E51').Value = '  +0.74%  '
